$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3   = "ppe || 18_hazards_to_humans_and_domestic_animals"
    4   = "application instructions || env warning - species"
    9   = "154_pesticide_storage"
    15  = "application instructions"
    16  = "application instructions"
    17  = "application instructions"
    18  = "application instructions"
    19  = "application instructions"
    20  = "application instructions"
    21  = "application instructions"
    22  = "application instructions"
    23  = "application instructions"
    24  = "application instructions"
    25  = "application instructions"
    26  = "application instructions"
    27  = "application instructions"
    28  = "application instructions"
    29  = "application instructions"
    30  = "application instructions"
    31  = "application instructions"
    32  = "application instructions"
    33  = "application instructions"
    34  = "application instructions"
    35  = "application instructions"
    36  = "application instructions"
    37  = "application instructions"
    38  = "application instructions"
    40  = "application instructions"
    41  = "application instructions"
    42  = "application instructions"
    43  = "application instructions"
    44  = "application instructions"
    45  = "application instructions"
    46  = "application instructions"
    47  = "application instructions"
    48  = "application instructions"
    49  = "application instructions"
    50  = "use restrictions"
    54  = "application instructions"
    55  = "application instructions"
    72  = "application instructions"
    153 = "application instructions"
    156 = "application instructions"
    167 = "application instructions"
    207 = "application instructions"
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
